$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '29.299.56'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.859.92'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '0.7036'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '238.25'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.07889'
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.3046'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '24.73'
$ws.Range('E10').Value = '  +6.33%  '
$ws.Range('D11').Value = '2.225.95'
$ws.Range('E11').Value = '  +19.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '0.08181'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '5.227'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '0.7177'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '89.65'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').Value = '29.295.85'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '5.831'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '0.000007801'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '13.24'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '238.60'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.074.33'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '7.565'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '162.85'
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '8.908'
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '0.1429'
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '18.10'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '1.920'
$ws.Range('E29').Value = '  -6.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '1.374'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '1.477'
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '4.055'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '0.05177'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '1.177'
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.7123'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '1.006'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.673'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.01852'
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '2.691'
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('D41').Value = '1.170.90'
$ws.Range('E41').Value = '  +2.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.9235'
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '5.983'
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '70.92'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.4258'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '101.66'
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.5332'
$ws.Range('E48').Value = '  -2.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.755'
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '9.180'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '1.965.54'
$ws.Range('E51').Value = '  -2.17%  '
